$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1: rows 2-11 (B=Insertion Sort, C=Selection Sort, D=Shell Sort) ---

$ws.Range("B2").Value = 796.88
$ws.Range("C2").Value = 578.13
$ws.Range("D2").Value = 31.79

$ws.Range("B3").Value = 3234.38
$ws.Range("C3").Value = 2265.63
$ws.Range("D3").Value = 62.5

$ws.Range("B4").Value = 12828.13
$ws.Range("C4").Value = 9125.69
$ws.Range("D4").Value = 156.33000000000001

$ws.Range("B5").Value = 51750.05
$ws.Range("C5").Value = 38125.06
$ws.Range("D5").Value = 343.44

$ws.Range("B6").Value = 214515.22
$ws.Range("C6").Value = 1688593.75
$ws.Range("D6").Value = 812.53

$ws.Range("B7").Value = 943625.99
$ws.Range("C7").Value = 724437.5
$ws.Range("D7").Value = 1796.21

$ws.Range("B8").ClearContents() | Out-Null
$ws.Range("C8").ClearContents() | Out-Null
$ws.Range("D8").Value = 4390.53

$ws.Range("B9").ClearContents() | Out-Null
$ws.Range("C9").ClearContents() | Out-Null
$ws.Range("D9").Value = 10750.22

$ws.Range("B10").ClearContents() | Out-Null
$ws.Range("C10").ClearContents() | Out-Null
$ws.Range("D10").Value = 26250.16

$ws.Range("B11").ClearContents() | Out-Null
$ws.Range("C11").ClearContents() | Out-Null
$ws.Range("D11").ClearContents() | Out-Null

# --- Table 2: rows 15-24 ---

$ws.Range("B15").Value = 30734.38
$ws.Range("C15").Value = 30640.91
$ws.Range("D15").Value = 1848.96

$ws.Range("B16").Value = 258234.3
$ws.Range("C16").Value = 280906.3
$ws.Range("D16").Value = 7781.25

$ws.Range("B17").ClearContents() | Out-Null
$ws.Range("C17").ClearContents() | Out-Null
$ws.Range("D17").Value = 37382.81

$ws.Range("B18").ClearContents() | Out-Null
$ws.Range("C18").ClearContents() | Out-Null
$ws.Range("D18").Value = 232937.5

$ws.Range("B19").ClearContents() | Out-Null
$ws.Range("C19").ClearContents() | Out-Null
$ws.Range("D19").ClearContents() | Out-Null

$ws.Range("B20").ClearContents() | Out-Null
$ws.Range("C20").ClearContents() | Out-Null
$ws.Range("D20").ClearContents() | Out-Null

$ws.Range("B21").ClearContents() | Out-Null
$ws.Range("C21").ClearContents() | Out-Null
$ws.Range("D21").ClearContents() | Out-Null

$ws.Range("B22").ClearContents() | Out-Null
$ws.Range("C22").ClearContents() | Out-Null
$ws.Range("D22").ClearContents() | Out-Null

$ws.Range("B23").ClearContents() | Out-Null
$ws.Range("C23").ClearContents() | Out-Null
$ws.Range("D23").ClearContents() | Out-Null

$ws.Range("B24").ClearContents() | Out-Null
$ws.Range("C24").ClearContents() | Out-Null
$ws.Range("D24").ClearContents() | Out-Null

# --- Re-format the freshly re-typed numbers: plain "General" number format,
#     keeping the table's "Dax-Regular" body font (this reproduces the
#     workbook's new cellXfs entry: numFmtId=0, fontId=1, applyFont=1).
#     Seed the style on B2 by touching its font directly, then fan it out to
#     every other re-typed cell with a formats-only paste so we don't mint a
#     duplicate style per cell. ---

$ws.Range("B2").Font.Name = "Dax-Regular"
$ws.Range("B2").Copy() | Out-Null

$formattedCells = "C2,D2,B3,C3,D3,B4,C4,D4,B5,C5,D5,B6,C6,D6,B7,C7,D7,D8,D9,D10,B15,C15,D15,B16,C16,D16,D17,D18"
foreach ($addr in $formattedCells.Split(",")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Restore the active-cell selection left by the author ---
$ws.Range("G19").Select() | Out-Null
